$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-21 Friday" "2025-02-22 Saturday"

Replace-Text "319÷9=35, 4" "217÷7=31, 0"
Replace-Text "487÷9=54, 1" "236÷9=26, 2"
Replace-Text "516÷5=103, 1" "839÷4=209, 3"
Replace-Text "698÷2=349, 0" "838÷6=139, 4"
Replace-Text "200÷7=28, 4" "299÷5=59, 4"

Replace-Text "836÷9=92, 8" "632÷8=79, 0"
Replace-Text "360÷9=40, 0" "540÷3=180, 0"
Replace-Text "694÷8=86, 6" "103÷7=14, 5"
Replace-Text "711÷8=88, 7" "467÷7=66, 5"
Replace-Text "417÷5=83, 2" "170÷5=34, 0"

Replace-Text "390÷9=43, 3" "830÷9=92, 2"
Replace-Text "382÷2=191, 0" "845÷3=281, 2"
Replace-Text "331÷2=165, 1" "410÷9=45, 5"
Replace-Text "571÷5=114, 1" "575÷2=287, 1"
Replace-Text "595÷9=66, 1" "352÷2=176, 0"

Replace-Text "131÷2=65, 1" "753÷9=83, 6"
Replace-Text "433÷7=61, 6" "417÷2=208, 1"
Replace-Text "676÷8=84, 4" "511÷9=56, 7"
Replace-Text "566÷7=80, 6" "767÷9=85, 2"
Replace-Text "725÷9=80, 5" "318÷9=35, 3"

Replace-Text "203÷4=50, 3" "820÷6=136, 4"
Replace-Text "287÷6=47, 5" "803÷7=114, 5"
Replace-Text "616÷8=77, 0" "916÷9=101, 7"
Replace-Text "482÷2=241, 0" "322÷2=161, 0"
Replace-Text "657÷7=93, 6" "368÷6=61, 2"

$d.Save()
